$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "addstudent dropdown" source list in column E, alongside
# the existing Classroom data that lives in columns A:D.
$ws.Range("E1").Value = "addstudent dropdown"
$ws.Range("E2").Value = "Grade 1"
$ws.Range("E3").Value = "Grade 2"
$ws.Range("E4").Value = "Grade 3"
$ws.Range("E5").Value = 8
$ws.Range("E6").Value = "LKG"
$ws.Range("E7").Value = "UKG"

# Size column E to fit its new content (target width is 21 characters)
# and leave the selection on the last cell that was filled in, matching
# the saved workbook's view state.
$ws.Columns("E").ColumnWidth = 20.166666666666668
[void]$ws.Range("E7").Select()
